$d = $word.ActiveDocument
$p = $d.Paragraphs.Last
$r = $d.Range($p.Range.Start, $p.Range.Start)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>We will now put this all together in a simulation.</w:t></w:r><w:r><w:t xml:space="preserve"> Lay out the index cards with some writing implements. Have each student (and each instructor except yourself) to run the first code block. It </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>give</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> them directions for which</w:t></w:r><w:r><w:t xml:space="preserve"> color</w:t></w:r><w:r><w:t xml:space="preserve"> to choose and tell them to write </w:t></w:r><w:r><w:t>number</w:t></w:r><w:r><w:t xml:space="preserve"> on the front</w:t></w:r><w:r><w:t>, and</w:t></w:r><w:r><w:t xml:space="preserve"> a</w:t></w:r><w:r><w:t xml:space="preserve"> 0 or 1</w:t></w:r><w:r><w:t xml:space="preserve"> on the back. We will now do a walkthrough of a </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>simulation. Each person represents a patient, with the color of their index card representing</w:t></w:r><w:r><w:t xml:space="preserve"> their injury severity, the number of the front representing the time they arrive at the hospital, and the 0 or 1 on the back representing whether they survive if they turned away.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>The walkthrough:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>First set aside some number of desks or chairs representing the hospital beds. I typically try to aim for roughly 1/3 the number of beds vs the number of people with cards.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Have all the people with cards line up in order. Namely have them line up so that people with lower numbers are at the front of the line and people with larger numbers are at the back of the line.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Now walk through the simulation doing FCFS first. For example, for FCFS, assign people in line to beds until there are not more free beds and assign everyone else to the other side of the room. Then have the everyone with a 1 raise their hands and count how many survivors there are (including those in beds).</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Repeat this process for ORANGE ONLY. Walking through the line only assigning the ORANGES to beds until there are no more free beds. Repeat the process of counting survivors and compare with the FCFS policy. Sometime there will be free beds when using ORANGE ONLY. If this is the case make sure to point it out.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Now explain that this was just one simulation and in order to truly compare the policies we would want to walk through this simulation many, many times. This can be accomplished by running simulations on the computer. The second code black will do just that. Allow the students to play with the parameters in this code block but beware if they set a number too large of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/></w:rPr><w:t>num_reps</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> then it may take a long (OR REALLY LONG) time to run.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Discuss which policy works better and instruct the students to see if they can find different combinations of parameters where each policy is better.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>The final code black isn’t too important but will display images of the output of the simulations run above.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)
